$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header cell F1 with same style as neighboring header cells
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# time_taken values for each data row (rows 2-72)
$timeTaken = @(
    "2021-10-05 10:50:49.775307",
    "2021-10-05 10:50:49.775318",
    "2021-10-05 10:50:49.775322",
    "2021-10-05 10:50:49.775324",
    "2021-10-05 10:50:49.775327",
    "2021-10-05 10:50:49.775330",
    "2021-10-05 10:50:49.775332",
    "2021-10-05 10:50:49.775335",
    "2021-10-05 10:50:49.775337",
    "2021-10-05 10:50:49.775340",
    "2021-10-05 10:50:49.775343",
    "2021-10-05 10:50:49.775345",
    "2021-10-05 10:50:49.775348",
    "2021-10-05 10:50:49.775351",
    "2021-10-05 10:50:49.775353",
    "2021-10-05 10:50:49.775356",
    "2021-10-05 10:50:49.775358",
    "2021-10-05 10:50:49.775361",
    "2021-10-05 10:50:49.775364",
    "2021-10-05 10:50:49.775366",
    "2021-10-05 10:50:49.775369",
    "2021-10-05 10:50:49.775371",
    "2021-10-05 10:50:49.775374",
    "2021-10-05 10:50:49.775376",
    "2021-10-05 10:50:49.775379",
    "2021-10-05 10:50:49.775382",
    "2021-10-05 10:50:49.775384",
    "2021-10-05 10:50:49.775387",
    "2021-10-05 10:50:49.775389",
    "2021-10-05 10:50:49.775392",
    "2021-10-05 10:50:49.775394",
    "2021-10-05 10:50:49.775397",
    "2021-10-05 10:50:49.775400",
    "2021-10-05 10:50:49.775402",
    "2021-10-05 10:50:49.775405",
    "2021-10-05 10:50:49.775407",
    "2021-10-05 10:50:49.775409",
    "2021-10-05 10:50:49.775412",
    "2021-10-05 10:50:49.775414",
    "2021-10-05 10:50:49.775417",
    "2021-10-05 10:50:49.775420",
    "2021-10-05 10:50:49.775422",
    "2021-10-05 10:50:49.775425",
    "2021-10-05 10:50:49.775427",
    "2021-10-05 10:50:49.775430",
    "2021-10-05 10:50:49.775432",
    "2021-10-05 10:50:49.775435",
    "2021-10-05 10:50:49.775437",
    "2021-10-05 10:50:49.775440",
    "2021-10-05 10:50:49.775442",
    "2021-10-05 10:50:49.775445",
    "2021-10-05 10:50:49.775447",
    "2021-10-05 10:50:49.775450",
    "2021-10-05 10:50:49.775453",
    "2021-10-05 10:50:49.775455",
    "2021-10-05 10:50:49.775458",
    "2021-10-05 10:50:49.775460",
    "2021-10-05 10:50:49.775463",
    "2021-10-05 10:50:49.775465",
    "2021-10-05 10:50:49.775468",
    "2021-10-05 10:50:49.775470",
    "2021-10-05 10:50:49.775473",
    "2021-10-05 10:50:49.775475",
    "2021-10-05 10:50:49.775478",
    "2021-10-05 10:50:49.775481",
    "2021-10-05 10:50:49.775484",
    "2021-10-05 10:50:49.775486",
    "2021-10-05 10:50:49.775489",
    "2021-10-05 10:50:49.775491",
    "2021-10-05 10:50:49.775493",
    "2021-10-05 10:50:49.775496"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

Write-Host "Done"
